# ISIS1225 - Lab 7 "Entrega Final" edit
# Updates the measured PROBING / CHAINING results with the final run's
# numbers, turns the old "+40/+30 demo" formulas into plain recorded
# values, fills in the previously-empty row of the CHAINING table, and
# adds the extra "Column1" column that Excel appends when the table is
# resized one column wider.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos Lab7")

# ---------------------------------------------------------------
# Tabla "Carga de Catálogo PROBING" (Table1, A2:C6)
# ---------------------------------------------------------------
$ws.Range("B3").Value = 25252.5
$ws.Range("C3").Value = 103.383

$ws.Range("B4").Value = 25235.3
$ws.Range("C4").Value = 346.592

$ws.Range("B5").Value = 25235.200000000001
$ws.Range("C5").Value = 247.185

$ws.Range("B6").Value = 25235.3
$ws.Range("C6").Value = 257.388

# ---------------------------------------------------------------
# Tabla "Carga de Catálogo CHAINING" (Table13, A10:C14 -> A10:D14)
# ---------------------------------------------------------------
$ws.Range("B11").Value = 25240.6
$ws.Range("C11").Value = 106.062

$ws.Range("B12").Value = 25235.200000000001
$ws.Range("C12").Value = 118.133

# Row 13 used to be left blank; the final run fills it in too.
$ws.Range("B13").Value = 25235.200000000001
$ws.Range("C13").Value = 122.469
$ws.Range("B12:C12").Copy() | Out-Null
$ws.Range("B13:C13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("B14").Value = 25235.200000000001
$ws.Range("C14").Value = 108.911

# Widen the CHAINING table by one column - Excel names the new column
# "Column1" and the table / autofilter range grows to A10:D14.
$lo = $ws.ListObjects.Item("Table13")
$lo.Resize($ws.Range("A10:D14"))
$ws.Range("D10").Value = "Column1"

# The new column's empty data cells pick up the same number format /
# centered+wrapped alignment as the rest of the table body.
$dataCells = $ws.Range("D11:D14")
$dataCells.NumberFormat = "0.00"
$dataCells.HorizontalAlignment = -4108
$dataCells.VerticalAlignment = -4108
$dataCells.WrapText = $true

# Final selection left on the sheet.
$ws.Range("D16").Select()
